$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 78640.71000000001
$ws.Range("I9").Value = 91703.414
$ws.Range("J9").Value = 264.5
$ws.Range("K9").Value = 91703.414
$ws.Range("L9").Value = 264.5
$ws.Range("M9").Value = -91534.414
$ws.Range("N9").Value = -602.5
$ws.Range("H32").Value = 3476.8
$ws.Range("I32").Value = 3472
$ws.Range("K32").Value = 3472
$ws.Range("M32").Value = -3146
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H112").Value = 2169.75
$ws.Range("J112").Value = 2339.7693
$ws.Range("L112").Value = 7019.3079
$ws.Range("N112").Value = -9235.3079
$ws.Range("H118").Value = 3057.75
$ws.Range("I118").Value = 2780.2856
$ws.Range("K118").Value = 8340.856800000001
$ws.Range("M118").Value = -6683.856800000001
$ws.Range("H137").Value = 2409.375
$ws.Range("J137").Value = 2330
$ws.Range("L137").Value = 6990
$ws.Range("N137").Value = -12090
$ws.Range("H138").Value = 7907.244
$ws.Range("J138").Value = 6573.647
$ws.Range("L138").Value = 19720.941
$ws.Range("N138").Value = -30000.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 450
$ws.Range("I4").Value = 450
$ws.Range("K4").Value = 450
$ws.Range("M4").Value = -334
$ws.Range("H5").Value = 423
$ws.Range("I5").Value = 423
$ws.Range("K5").Value = 423
$ws.Range("M5").Value = -311
$ws.Range("H32").Value = 23227.215
$ws.Range("I32").Value = 18365.084
$ws.Range("K32").Value = 18365.084
$ws.Range("M32").Value = -18078.084
$ws.Range("H45").Value = 2501.3572
$ws.Range("I45").Value = 1121.4
$ws.Range("J45").Value = 3268
$ws.Range("K45").Value = 1121.4
$ws.Range("L45").Value = 3268
$ws.Range("M45").Value = -744.4000000000001
$ws.Range("N45").Value = -4022

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 423
$ws.Range("I4").Value = 423
$ws.Range("K4").Value = 423
$ws.Range("M4").Value = -308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 741.0833
$ws.Range("J22").Value = 992.8
$ws.Range("L22").Value = 992.8
$ws.Range("N22").Value = -1692.8
$ws.Range("H31").Value = 66519.25
$ws.Range("I31").Value = 4192.143
$ws.Range("J31").Value = 114995.89
$ws.Range("K31").Value = 4192.143
$ws.Range("L31").Value = 114995.89
$ws.Range("M31").Value = -3897.143
$ws.Range("N31").Value = -115585.89
$ws.Range("H34").Value = 66519.25
$ws.Range("I34").Value = 4192.143
$ws.Range("J34").Value = 114995.89
$ws.Range("K34").Value = 4192.143
$ws.Range("L34").Value = 114995.89
$ws.Range("M34").Value = -3990.143
$ws.Range("N34").Value = -115399.89
$ws.Range("H58").Value = 5178.533
$ws.Range("I58").Value = 4746.0835
$ws.Range("J58").Value = 6908.3335
$ws.Range("K58").Value = 4746.0835
$ws.Range("L58").Value = 6908.3335
$ws.Range("M58").Value = -4543.0835
$ws.Range("N58").Value = -7314.3335
$ws.Range("H86").Value = 6171.636
$ws.Range("I86").Value = 4986.125
$ws.Range("J86").Value = 9333
$ws.Range("K86").Value = 4986.125
$ws.Range("L86").Value = 9333
$ws.Range("M86").Value = -3863.125
$ws.Range("N86").Value = -11579
$ws.Range("H89").Value = 6171.636
$ws.Range("I89").Value = 4986.125
$ws.Range("J89").Value = 9333
$ws.Range("K89").Value = 24930.625
$ws.Range("L89").Value = 46665
$ws.Range("M89").Value = -19314.625
$ws.Range("N89").Value = -57897
$ws.Range("H105").Value = 975.25
$ws.Range("I105").Value = 828.8570999999999
$ws.Range("K105").Value = 828.8570999999999
$ws.Range("M105").Value = 918.1429000000001
$ws.Range("H136").Value = 5178.533
$ws.Range("I136").Value = 4746.0835
$ws.Range("J136").Value = 6908.3335
$ws.Range("K136").Value = 14238.2505
$ws.Range("L136").Value = 20725.0005
$ws.Range("M136").Value = -11688.2505
$ws.Range("N136").Value = -25825.0005
$ws.Range("H141").Value = 493607
$ws.Range("J141").Value = 578398.7
$ws.Range("L141").Value = 578398.7
$ws.Range("N141").Value = -588758.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1883391.9
$ws.Range("I11").Value = 1883391.9
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5650175.699999999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -5650035.699999999
$ws.Range("N11").ClearContents()
$ws.Range("H32").Value = 750241.5
$ws.Range("I32").Value = 666989.3
$ws.Range("J32").Value = 999998
$ws.Range("K32").Value = 2000967.9
$ws.Range("L32").Value = 2999994
$ws.Range("M32").Value = -2000684.9
$ws.Range("N32").Value = -3000560
$ws.Range("H69").Value = 866.3333
$ws.Range("I69").Value = 866.3333
$ws.Range("K69").Value = 2598.9999
$ws.Range("M69").Value = -1787.9999
$ws.Range("H72").Value = 866.3333
$ws.Range("I72").Value = 866.3333
$ws.Range("K72").Value = 7796.9997
$ws.Range("M72").Value = -3740.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2999.5
$ws.Range("I80").Value = 2999.5
$ws.Range("K80").Value = 2999.5
$ws.Range("M80").Value = -2001.5
$ws.Range("H83").Value = 2999.5
$ws.Range("I83").Value = 2999.5
$ws.Range("K83").Value = 14997.5
$ws.Range("M83").Value = -10005.5
$ws.Range("H102").Value = 16305.789
$ws.Range("I102").Value = 18816
$ws.Range("J102").Value = 2918
$ws.Range("K102").Value = 18816
$ws.Range("L102").Value = 2918
$ws.Range("M102").Value = -17194
$ws.Range("N102").Value = -6162

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1056.174
$ws.Range("J22").Value = 1368.375
$ws.Range("L22").Value = 1368.375
$ws.Range("N22").Value = -1958.375
$ws.Range("H27").Value = 1056.174
$ws.Range("J27").Value = 1368.375
$ws.Range("L27").Value = 1368.375
$ws.Range("N27").Value = -1582.375
$ws.Range("H40").Value = 3425.3635
$ws.Range("I40").Value = 3554.4285
$ws.Range("J40").Value = 3199.5
$ws.Range("K40").Value = 3554.4285
$ws.Range("L40").Value = 3199.5
$ws.Range("M40").Value = -3418.4285
$ws.Range("N40").Value = -3471.5
$ws.Range("H55").Value = 979.86664
$ws.Range("I55").Value = 821.7
$ws.Range("K55").Value = 821.7
$ws.Range("M55").Value = -648.7
$ws.Range("H68").Value = 2040.6923
$ws.Range("I68").Value = 1957.4546
$ws.Range("K68").Value = 1957.4546
$ws.Range("M68").Value = -1208.4546
$ws.Range("H71").Value = 2040.6923
$ws.Range("I71").Value = 1957.4546
$ws.Range("K71").Value = 9787.273000000001
$ws.Range("M71").Value = -6043.273000000001
$ws.Range("H82").Value = 4715.364
$ws.Range("I82").Value = 921.25
$ws.Range("K82").Value = 921.25
$ws.Range("M82").Value = -560.25
$ws.Range("H85").Value = 4715.364
$ws.Range("I85").Value = 921.25
$ws.Range("K85").Value = 921.25
$ws.Range("M85").Value = 326.75
$ws.Range("H136").Value = 3576.375
$ws.Range("I136").Value = 3501.8333
$ws.Range("K136").Value = 10505.4999
$ws.Range("M136").Value = -7955.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 260000
$ws.Range("I5").Value = 20000
$ws.Range("K5").Value = 20000
$ws.Range("M5").Value = -19888
$ws.Range("H62").Value = 4939.4
$ws.Range("I62").Value = 3924.5
$ws.Range("K62").Value = 3924.5
$ws.Range("M62").Value = -3300.5
$ws.Range("H65").Value = 4939.4
$ws.Range("I65").Value = 3924.5
$ws.Range("K65").Value = 19622.5
$ws.Range("M65").Value = -16502.5
$ws.Range("H96").Value = 1027.5714
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1027.5714
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1027.5714
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -3773.5714
$ws.Range("H99").Value = 50499.5
$ws.Range("J99").Value = 50499.5
$ws.Range("L99").Value = 50499.5
$ws.Range("N99").Value = -56489.5
$ws.Range("H100").Value = 1013.125
$ws.Range("I100").Value = 924
$ws.Range("J100").Value = 1161.6666
$ws.Range("K100").Value = 1848
$ws.Range("L100").Value = 2323.3332
$ws.Range("M100").Value = -1307
$ws.Range("N100").Value = -3405.3332
$ws.Range("H107").Value = 3281.4546
$ws.Range("I107").Value = 3067.3333
$ws.Range("J107").Value = 3361.75
$ws.Range("K107").Value = 9201.999899999999
$ws.Range("L107").Value = 10085.25
$ws.Range("M107").Value = -7281.999899999999
$ws.Range("N107").Value = -13925.25
$ws.Range("H132").Value = 51744.5
$ws.Range("I132").Value = 99994
$ws.Range("K132").Value = 299982
$ws.Range("M132").Value = -297452
$ws.Range("H136").Value = 10459.176
$ws.Range("I136").Value = 10371.881
$ws.Range("J136").Value = 10703.6
$ws.Range("K136").Value = 31115.643
$ws.Range("L136").Value = 32110.8
$ws.Range("M136").Value = -28565.643
$ws.Range("N136").Value = -37210.8

Write-Host "Applied all updates"